# The deck's Design ("Integral") theme is switched back to the stock
# "Office Theme" colour palette (Design > Variants > Colors > Office in the
# PowerPoint UI). Only the theme colour scheme (a:clrScheme, 12 slots)
# differs between the "Integral" palette and the default "Office" palette -
# the font scheme and format scheme are identical - so the edit is applied
# through the documented ThemeColorScheme object on the slide master's
# Theme, one RGB slot at a time (RGB values are COM's 0x00BBGGRR longs).

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0          # dk1      -> 000000
$tcs.Colors(2).RGB  = 16777215   # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      -> 44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  -> FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  -> 70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    -> 0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink -> 954F72
